$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update AdductLabels default to include sodium and potassium in addition to proton
$ws.Range("B8").Value = "proton,sodium,potassium"

# Update AdductMasses default to hold the matching comma-separated masses (now text, not a single number)
$ws.Range("B9").Value = "1.00727647,22.989769,39.0983"

# Column B needs to be widened to fit the new longer default values
$ws.Columns.Item(2).ColumnWidth = 35.92

# Reflect the final selection location left in the saved workbook
$ws.Range("B5").Select()
